# "Edited Workplan - Day 4"
# Updates the "Work Plan" sheet's Day 4 section (rows 26-30):
#  - Row 26 (Deploy App from Ansible to Tomcat): clears Actual End Date, drops percent
#    complete to 86 and flips status from Completed to In Progress.
#  - Row 27 (Test app using Selenium): records an Actual Start Date, sets percent
#    complete to 70 and flips status to In Progress.
#  - Row 29 (Debug Project Pipeline): records an Actual Start Date, sets percent
#    complete to 70 and flips status to In Progress.
#  - Row 30 (Finalize Documents): records an Actual Start Date, sets percent
#    complete to 62 and flips status to In Progress.
# Also moves the sheet's cursor/selection to B26 and clears the frozen top-left cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Plan")

# ---- Row 26: Deploy App from Ansible to Tomcat ----
$ws.Range("F26").ClearContents()
$ws.Range("G26").Value = 86
$ws.Range("H26").Value = "In Progress"

# ---- Row 27: Test app using Selenium ----
$ws.Range("E27").Value = 42569
$ws.Range("G27").Value = 70

$ws.Range("A8").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = "In Progress"

# ---- Row 29: Debug Project Pipeline ----
$ws.Range("C27").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = 42569

$ws.Range("G22").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 70

$ws.Range("A8").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = "In Progress"

# ---- Row 30: Finalize Documents ----
$ws.Range("C27").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = 42570

$ws.Range("G30").Value = 62
$ws.Range("G30").HorizontalAlignment = -4108

$ws.Range("A8").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = "In Progress"

$excel.CutCopyMode = $false

# ---- Selection / view ----
$ws.Range("B26").Select()
